$wb = $excel.ActiveWorkbook

# Update the CompanyName value on the "Contact" sheet (row 2 / column A)
# from "StandardTestCompany" to "ActivityCompany".
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("A2").Value = "ActivityCompany"

# Make "Contact" the active sheet and set its selection to B10.
# (This also clears the "tabSelected" flag on the previously active
# "Users" sheet, whose own selection (C5) is left untouched.)
$wsContact.Activate()
$wsContact.Range("B10").Select()
